# edit.ps1
# Commit: "Change cooling and PV"
#   - Cooling column Y values all multiplied by COP = 4.3 (and downstream
#     cash/energy columns recomputed accordingly)
#   - Re-run with new PV profiles, adjusted to meet total production kWh
#   - (Failed try changing HW storage settings -- left as-is here, the
#     numeric results below are the final, kept state)
#
# This script reproduces the data + light formatting changes made to
# Sheet1 of the results database workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update numeric result cells (recomputed cooling * COP=4.3, new PV profiles) ---
# Row 2
$ws.Range("D2").Value = 23484.81
$ws.Range("E2").Value = 1992.98
$ws.Range("F2").Value = 6910
$ws.Range("G2").Value = 1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 4132.87
$ws.Range("S2").Value = 1645.58
$ws.Range("T2").Value = 2908.77
$ws.Range("U2").Value = 404.46
$ws.Range("V2").Value = 1628.56
$ws.Range("W2").Value = 581.75
$ws.Range("X2").Value = -18.850000000000001
$ws.Range("Y2").Value = 837.91
$ws.Range("Z2").Value = 905.07
$ws.Range("AA2").Value = 1075.07
$ws.Range("AB2").Value = 41206.579299999998
# Row 3
$ws.Range("D3").Value = 18151.439999999999
$ws.Range("E3").Value = 1364.91
$ws.Range("F3").Value = 6800
$ws.Range("H3").Value = 2
$ws.Range("R3").Value = 1935.67
$ws.Range("S3").Value = 941.42
$ws.Range("T3").Value = 1935.67
$ws.Range("W3").Value = 387.13
$ws.Range("Y3").Value = 637.13
$ws.Range("Z3").Value = 517.78
$ws.Range("AA3").Value = 687.78
$ws.Range("AB3").Value = 25287.2893
# Row 4
$ws.Range("D4").Value = 26066.29
$ws.Range("E4").Value = 2034.04
$ws.Range("R4").Value = 3982.3
$ws.Range("S4").Value = 1368.32
$ws.Range("T4").Value = 3982.3
$ws.Range("W4").Value = 796.46
$ws.Range("Y4").Value = 1071.46
$ws.Range("Z4").Value = 752.58
$ws.Range("AA4").Value = 922.58
$ws.Range("AB4").Value = 39110.735000000001
# Row 5
$ws.Range("D5").Value = 30650.97
$ws.Range("E5").Value = 2302.7399999999998
$ws.Range("R5").Value = 4168.38
$ws.Range("S5").Value = 1743.75
$ws.Range("T5").Value = 4168.38
$ws.Range("W5").Value = 833.68
$ws.Range("Y5").Value = 1133.68
$ws.Range("Z5").Value = 959.06
$ws.Range("AA5").Value = 1129.06
$ws.Range("AB5").Value = 47489.95
# Row 6
$ws.Range("D6").Value = 24825.84
$ws.Range("E6").Value = 1942.6
$ws.Range("F6").Value = 8670
$ws.Range("G6").Value = 2
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 4
$ws.Range("R6").Value = 5167.28
$ws.Range("S6").Value = 1426.17
$ws.Range("T6").Value = 3475.46
$ws.Range("U6").Value = 1328.06
$ws.Range("V6").Value = 3019.87
$ws.Range("W6").Value = 695.09
$ws.Range("X6").Value = -61.89
$ws.Range("Y6").Value = 908.21
$ws.Range("Z6").Value = 784.4
$ws.Range("AA6").Value = 954.4
$ws.Range("AB6").Value = 39188.255599999997
# Row 7
$ws.Range("D7").Value = 21071.89
$ws.Range("E7").Value = 1728.1
$ws.Range("F7").Value = 6700
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("R7").Value = 1677.58
$ws.Range("S7").Value = 1741.05
$ws.Range("T7").Value = 1677.58
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 335.52
$ws.Range("X7").Value = 0
$ws.Range("Y7").Value = 560.52
$ws.Range("Z7").Value = 957.58
$ws.Range("AA7").Value = 1127.58
$ws.Range("AB7").Value = 42412.075100000002
# Row 8
$ws.Range("D8").Value = 31834.38
$ws.Range("E8").Value = 2727.6
$ws.Range("R8").Value = 4004.17
$ws.Range("S8").Value = 2621.39
$ws.Range("T8").Value = 4004.17
$ws.Range("W8").Value = 800.83
$ws.Range("Y8").Value = 1075.83
$ws.Range("Z8").Value = 1441.77
$ws.Range("AA8").Value = 1611.77
$ws.Range("AB8").Value = 68485.450700000001
# Row 9
$ws.Range("D9").Value = 38440.17
$ws.Range("E9").Value = 3239.32
$ws.Range("R9").Value = 4186.0600000000004
$ws.Range("S9").Value = 3440.2
$ws.Range("T9").Value = 4186.0600000000004
$ws.Range("W9").Value = 837.21
$ws.Range("Y9").Value = 1137.21
$ws.Range("Z9").Value = 1892.11
$ws.Range("AA9").Value = 2062.11
$ws.Range("AB9").Value = 85861.739400000006
# Row 10
$ws.Range("D10").Value = 25690.400000000001
$ws.Range("E10").Value = 1878.22
$ws.Range("F10").Value = 10070
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 3
$ws.Range("R10").Value = 2899.97
$ws.Range("S10").Value = 2204.3000000000002
$ws.Range("T10").Value = 827.81
$ws.Range("U10").Value = 852.11
$ws.Range("V10").Value = 2724.15
$ws.Range("W10").Value = 165.56
$ws.Range("X10").Value = -39.71
$ws.Range("Y10").Value = 375.85
$ws.Range("Z10").Value = 1212.3599999999999
$ws.Range("AA10").Value = 1382.36
$ws.Range("AB10").Value = 47453.356800000001
# Row 11
$ws.Range("D11").Value = 23074.720000000001
$ws.Range("E11").Value = 1968.92
$ws.Range("F11").Value = 6700
$ws.Range("H11").Value = 1
$ws.Range("R11").Value = 1290.45
$ws.Range("S11").Value = 2319.69
$ws.Range("T11").Value = 1290.45
$ws.Range("W11").Value = 258.08999999999997
$ws.Range("Y11").Value = 483.09
$ws.Range("Z11").Value = 1275.83
$ws.Range("AA11").Value = 1445.83
$ws.Range("AB11").Value = 53043.0864
# Row 12
$ws.Range("D12").Value = 35765.550000000003
$ws.Range("E12").Value = 3200.29
$ws.Range("R12").Value = 3912.73
$ws.Range("S12").Value = 3514.08
$ws.Range("T12").Value = 3912.73
$ws.Range("W12").Value = 782.55
$ws.Range("Y12").Value = 1057.55
$ws.Range("Z12").Value = 1932.75
$ws.Range("AA12").Value = 2102.75
$ws.Range("AB12").Value = 87381.915500000003
# Row 13
$ws.Range("D13").Value = 43584.84
$ws.Range("E13").Value = 3857.92
$ws.Range("R13").Value = 4108.26
$ws.Range("S13").Value = 4593.22
$ws.Range("T13").Value = 4108.26
$ws.Range("W13").Value = 821.65
$ws.Range("Y13").Value = 1121.6500000000001
$ws.Range("Z13").Value = 2526.27
$ws.Range("AA13").Value = 2696.27
$ws.Range("AB13").Value = 110050.3893

# --- Column B/C (city/size) lose their explicit centered/bordered style on
#     the data rows -- revert those cells back to the workbook's default
#     "Normal" style (drops the s="5" style index) ---
foreach ($r in 2..13) {
    $ws.Range("B$r").Style = "Normal"
    $ws.Range("C$r").Style = "Normal"
}

# --- Column A (house) data cells (A2:A12) lose their bottom border edge ---
foreach ($r in 2..12) {
    $ws.Range("A$r").Borders(9).LineStyle = -4142
}

# --- Selection moves from F26 to G26 ---
$ws.Range("G26").Select()
